$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update: force Text number format first so Excel
# does not auto-convert numeric-looking strings (e.g. '96.15') into
# floating point numbers, which would not match the source inlineStr text.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.637.16'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.287.32'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '96.15'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '267.69'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.27%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.81'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0934'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.89'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.54%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.630.14'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.13'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.849'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.289.27'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.581.52'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.67%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.21'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.51'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +10.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.69'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.12'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.43%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.20'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.16'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.22%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.20'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.83'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.23%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.36'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.31%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.17%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.35'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.42'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.242'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.71%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.24'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.35'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.90'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.79'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.20'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.05%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '97.37'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.00%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.508.61'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.184'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +8.44%  '
